$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Albedo Generator" table used to end with V1.03.2 (row 9) and V1.03.6
# (row 10), followed by a gap and the "Dependencies here" note at row 14.
# The albedo extractor's network architecture changed to an FFA-net style,
# which produced new benchmark entries: V1.03.4, V1.0.3.7, V1.04.1 and
# V1.04.2. Re-create the table with the new rows inserted in version order,
# shifting the trailing "Dependencies here" note down accordingly.

# Insert a new row at 10, pushing the existing "V1.03.6" row (currently row
# 10) down to row 11, and the "Dependencies here" row (currently row 14)
# down to row 15.
$ws.Rows("10:10").Insert()

# New row 10: V1.03.4 results.
$ws.Range("A10").Value = "V1.03.4"
$ws.Range("B10").Value = 16.04551
$ws.Range("C10").Value = 0.02487
$ws.Range("D10").Value = 0.6146

# Rows 12-14 are still empty spacer rows at this point (row 11 holds the
# "V1.03.6" entry that was shifted down). Fill them in version order so the
# newly-created shared-string entries line up with the source order
# (V1.03.4, V1.04.1, V1.04.2, V1.0.3.7).
$ws.Range("A13").Value = "V1.04.1"
$ws.Range("B13").Value = 18.64907
$ws.Range("C13").Value = 0.01366
$ws.Range("D13").Value = 0.66189

$ws.Range("A14").Value = "V1.04.2"
$ws.Range("B14").Value = 17.95273
$ws.Range("C14").Value = 0.01603
$ws.Range("D14").Value = 0.63198

$ws.Range("A12").Value = "V1.0.3.7"
$ws.Range("B12").Value = 18.28864
$ws.Range("C12").Value = 0.01485
$ws.Range("D12").Value = 0.70436

# Insert 2 more blank rows before the "Dependencies here" row (now at row
# 15) to push it down to row 17, keeping the same visual gap as before.
$ws.Rows("15:16").Insert()

# Move the threaded "Dependencies" comment from the old note location (A14)
# to its new location (A17).
$oldComment = $ws.Range("A14").Comment
$commentText = $oldComment.Text()
$ws.Range("A17").AddCommentThreaded($commentText) | Out-Null
$oldComment.Delete()

# Update the active selection to reflect where editing left off.
$ws.Range("E15").Select()
